$d = $word.ActiveDocument

# --- Update the title date line -------------------------------------------
$d.Paragraphs.Item(1).Range.Text = "2025-07-18 Friday"

# --- Update the practice-problem table --------------------------------------
$t = $d.Tables.Item(1)

# Map of 1-based table row -> new cell values (left to right).
# Only the rows that contain answer text are touched; the blank spacer
# rows in between are left untouched.
$updates = @{
    1  = @("42÷9=4, 6", "28÷8=3, 4", "77÷7=11, 0", "42÷6=7, 0", "99÷8=12, 3")
    5  = @("94÷8=11, 6", "34÷4=8, 2", "87÷8=10, 7", "95÷9=10, 5", "39÷6=6, 3")
    9  = @("60÷5=12, 0", "75÷2=37, 1", "74÷7=10, 4", "91÷3=30, 1", "12÷6=2, 0")
    13 = @("63÷2=31, 1", "59÷3=19, 2", "53÷2=26, 1", "84÷5=16, 4", "18÷7=2, 4")
    17 = @("11÷5=2, 1", "84÷9=9, 3", "34÷5=6, 4", "34÷4=8, 2", "56÷7=8, 0")
}

foreach ($rowIndex in $updates.Keys) {
    $row = $t.Rows.Item($rowIndex)
    $values = $updates[$rowIndex]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $cell = $row.Cells.Item($i + 1)
        $cell.Range.Text = $values[$i]
    }
}
